$wb = $excel.ActiveWorkbook

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1199.8125
$ws.Range("I111").Value = 1045.9231
$ws.Range("J111").Value = 1866.6666
$ws.Range("K111").Value = 3137.7693
$ws.Range("L111").Value = 5599.9998
$ws.Range("M111").Value = -70.76929999999993
$ws.Range("N111").Value = -11733.9998

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2404.2432
$ws.Range("I132").Value = 2483.5151
$ws.Range("J132").Value = 1750.25
$ws.Range("K132").Value = 7450.5453
$ws.Range("L132").Value = 5250.75
$ws.Range("M132").Value = -4920.5453
$ws.Range("N132").Value = -10310.75

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 975.913
$ws.Range("I137").Value = 838
$ws.Range("K137").Value = 2514
$ws.Range("M137").Value = 36

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3966.6482
$ws.Range("I138").Value = 1891.9048
$ws.Range("J138").Value = 5286.9395
$ws.Range("K138").Value = 5675.7144
$ws.Range("L138").Value = 15860.8185
$ws.Range("M138").Value = -535.7143999999998
$ws.Range("N138").Value = -26140.8185

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3771.9412
$ws.Range("I2").Value = 2939.4614
$ws.Range("J2").Value = 6477.5
$ws.Range("K2").Value = 2939.4614
$ws.Range("L2").Value = 6477.5
$ws.Range("M2").Value = -2826.4614
$ws.Range("N2").Value = -6703.5

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2080.152
$ws.Range("I45").Value = 1829.9
$ws.Range("J45").Value = 2549.375
$ws.Range("K45").Value = 1829.9
$ws.Range("L45").Value = 2549.375
$ws.Range("M45").Value = -1452.9
$ws.Range("N45").Value = -3303.375

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1799.5294
$ws.Range("J74").Value = 1962.5
$ws.Range("L74").Value = 1962.5
$ws.Range("N74").Value = -3710.5

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1799.5294
$ws.Range("J77").Value = 1962.5
$ws.Range("L77").Value = 9812.5
$ws.Range("N77").Value = -18548.5

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3771.9412
$ws.Range("I116").Value = 2939.4614
$ws.Range("J116").Value = 6477.5
$ws.Range("K116").Value = 2939.4614
$ws.Range("L116").Value = 6477.5
$ws.Range("M116").Value = -645.4614000000001
$ws.Range("N116").Value = -11065.5

# ARM row 125
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()  # was -59840

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 38501380
$ws.Range("I132").Value = 71429220
$ws.Range("J132").Value = 85566
$ws.Range("K132").Value = 214287660
$ws.Range("L132").Value = 256698
$ws.Range("M132").Value = -214285130
$ws.Range("N132").Value = -261758

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3771.9412
$ws.Range("I3").Value = 2939.4614
$ws.Range("J3").Value = 6477.5
$ws.Range("K3").Value = 2939.4614
$ws.Range("L3").Value = 6477.5
$ws.Range("M3").Value = -2825.4614
$ws.Range("N3").Value = -6705.5

# BSM row 124
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6823.4414
$ws.Range("I134").Value = 2235.963
$ws.Range("J134").Value = 24518
$ws.Range("K134").Value = 6707.889000000001
$ws.Range("L134").Value = 73554
$ws.Range("M134").Value = -4172.889000000001
$ws.Range("N134").Value = -78624

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5321158
$ws.Range("I31").Value = 7354498
$ws.Range("J31").Value = 3192.3076
$ws.Range("K31").Value = 7354498
$ws.Range("L31").Value = 3192.3076
$ws.Range("M31").Value = -7354203
$ws.Range("N31").Value = -3782.3076

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5321158
$ws.Range("I34").Value = 7354498
$ws.Range("J34").Value = 3192.3076
$ws.Range("K34").Value = 7354498
$ws.Range("L34").Value = 3192.3076
$ws.Range("M34").Value = -7354296
$ws.Range("N34").Value = -3596.3076

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 711
$ws.Range("I5").Value = 660.3200000000001
$ws.Range("J5").Value = 1133.3334
$ws.Range("K5").Value = 1980.96
$ws.Range("L5").Value = 3400.0002
$ws.Range("M5").Value = -1868.96
$ws.Range("N5").Value = -3624.0002

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 448.81818
$ws.Range("I23").Value = 153.85715
$ws.Range("J23").Value = 965
$ws.Range("K23").Value = 461.57145
$ws.Range("L23").Value = 2895
$ws.Range("M23").Value = -226.57145
$ws.Range("N23").Value = -3365

# CUL row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 500
$ws.Range("I86").Value = 500
$ws.Range("K86").Value = 1500
$ws.Range("M86").Value = -314

# CUL row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 500
$ws.Range("I89").Value = 500
$ws.Range("K89").Value = 4500
$ws.Range("M89").Value = 1428

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1005.8933
$ws.Range("I113").Value = 833
$ws.Range("J113").Value = 1015.6338
$ws.Range("K113").Value = 2499
$ws.Range("L113").Value = 3046.9014
$ws.Range("M113").Value = -329
$ws.Range("N113").Value = -7386.9014

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 526.5517
$ws.Range("I122").Value = 306.83334
$ws.Range("J122").Value = 1581.2
$ws.Range("K122").Value = 2761.50006
$ws.Range("L122").Value = 14230.8
$ws.Range("M122").Value = -311.5000600000003
$ws.Range("N122").Value = -19130.8

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 12196301
$ws.Range("I131").Value = 2242.5
$ws.Range("J131").Value = 15152437
$ws.Range("K131").Value = 6727.5
$ws.Range("L131").Value = 45457311
$ws.Range("M131").Value = -1687.5
$ws.Range("N131").Value = -45467391

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 711
$ws.Range("I135").Value = 660.3200000000001
$ws.Range("J135").Value = 1133.3334
$ws.Range("K135").Value = 5942.88
$ws.Range("L135").Value = 10200.0006
$ws.Range("M135").Value = -3407.88
$ws.Range("N135").Value = -15270.0006

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 429.2857
$ws.Range("I107").Value = 320
$ws.Range("K107").Value = 320
$ws.Range("M107").Value = 1600

# GSM row 127
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 54980
$ws.Range("J127").Value = 54980
$ws.Range("L127").Value = 54980
$ws.Range("N127").Value = -64900

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 290871.25
$ws.Range("I132").Value = 38066.89
$ws.Range("J132").Value = 911391.0600000001
$ws.Range("K132").Value = 114200.67
$ws.Range("L132").Value = 2734173.18
$ws.Range("M132").Value = -111670.67
$ws.Range("N132").Value = -2739233.18

# GSM row 139
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()  # was -35280

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5055.44
$ws.Range("I46").Value = 1302.6
$ws.Range("J46").Value = 7557.3335
$ws.Range("K46").Value = 1302.6
$ws.Range("L46").Value = 7557.3335
$ws.Range("M46").Value = -1114.6
$ws.Range("N46").Value = -7933.3335

# LTW row 124
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820
